# Projektarbeit WS 18/19 - C_mu_50erWellen.xlsx
# "Abbildungen für Kapitel 6+ Änderungen an Kap 6"
#
# Updates the existing C_mu computation (new calibration constants),
# adds a second (0 rpm) computation in columns G/H, and adds a third
# block (rows 7-11) relating P_plenum readings to measured C_mu plus a
# standalone pressure -> velocity check in columns J/K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Existing table (A:C) - new physical constants in the formula
#    (written cell-by-cell so the pre-existing C4:C6 shared-formula
#    group keeps its si="0" id)
# ---------------------------------------------------------------
$ws.Range("C3").Formula = "=(A3*0.5*14.94^2*0.39*0.0534*1.17)/(2*0.384*0.19*10^(-3))"
$ws.Range("C4").Formula = "=(A4*0.5*14.94^2*0.39*0.0534*1.17)/(2*0.384*0.19*10^(-3))"
$ws.Range("C5").Formula = "=(A5*0.5*14.94^2*0.39*0.0534*1.17)/(2*0.384*0.19*10^(-3))"
$ws.Range("C6").Formula = "=(A6*0.5*14.94^2*0.39*0.0534*1.17)/(2*0.384*0.19*10^(-3))"

# ---------------------------------------------------------------
# 2) New table: "Werte für 0rpm" (G:H), rows 2-6
# ---------------------------------------------------------------
$ws.Range("G2").Value = "Werte für 0rpm"

$ws.Range("G3").Value = 0.064
$ws.Range("G4").Value = 0.27
$ws.Range("G5").Value = 0.416
$ws.Range("G6").Value = 0.582

$ws.Range("H3").Formula    = "=(G3*0.5*14.94^2*0.39*0.0534*1.17)/(2*0.384*0.19*10^(-3))"
$ws.Range("H4:H6").Formula = "=(G4*0.5*14.94^2*0.39*0.0534*1.17)/(2*0.384*0.19*10^(-3))"

# ---------------------------------------------------------------
# 3) New table: "Werte für rpm >0" / "pPlenum offen mittel" (G:K), rows 7-11
# ---------------------------------------------------------------
$ws.Range("G7").Value = "Werte für rpm >0"

# keep shared-string insertion order matching the source document:
# "Größerer Druck ..." (J9) is registered before "pPlenum offen mittel" (I7)
$ws.Range("J9").Value = "Größerer Druck nicht möglich, da sonst inkompressibel"
$ws.Range("I7").Value = "pPlenum offen mittel"

$ws.Range("H8").Value  = 1800
$ws.Range("H9").Value  = 4500
$ws.Range("H10").Value = 5400
$ws.Range("H11").Value = 6100

$ws.Range("G8").Formula     = "=(H8*2*0.384*0.00019)/(0.5*14.94^2*0.39*0.0534*1.17)"
$ws.Range("G9:G11").Formula = "=(H9*2*0.384*0.00019)/(0.5*14.94^2*0.39*0.0534*1.17)"

$ws.Range("I8").Value  = 1000
$ws.Range("I9").Value  = 2000
$ws.Range("I10").Value = 3000
$ws.Range("I11").Value = 4000

$ws.Range("J11").Value = 9000
$ws.Range("K11").Formula = "=SQRT(2*J11/1.17)"

# ---------------------------------------------------------------
# 4) Selection / window state, matching the saved UI state
# ---------------------------------------------------------------
$ws.Range("H3").Select()
$excel.WindowState = -4140
